# Apply the quarterly/hourly cryptocurrency price & volume refresh
# as produced by the GitHub Actions scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.635.20"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").Value = "1.534.27"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.36"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3944"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3163"
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.44"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07164"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.050"
$ws.Range("E11").Value = "  -6.66%  "
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.674"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.56"
$ws.Range("E14").Value = "  -4.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.597"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "1.572.79"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001088"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06608"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.66"
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  -4.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.43"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("E23").Value = "  -6.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.350"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "21.642.36"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.353"
$ws.Range("E26").Value = "  -7.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.78"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.844"
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("D30").Value = "1.715.89"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.07"
$ws.Range("E31").Value = "  -3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.985"
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9443"
$ws.Range("E33").Value = "  -15.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08165"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.509"
$ws.Range("E35").Value = "  -8.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.152"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06012"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02212"
$ws.Range("E38").Value = "  -3.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.449"
$ws.Range("E39").Value = "  -14.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2017"
$ws.Range("E40").Value = "  -4.04%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.177"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.94"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5757"
$ws.Range("E44").Value = "  -3.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.01"
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.709"
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5491"
$ws.Range("E47").Value = "  -4.55%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.876"
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "116.49"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  -2.94%  "
